$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.363.50'
$ws.Range("E2").Value = '  +1.34%  '
$ws.Range("D3").Value = '1.627.42'
$ws.Range("E3").Value = '  +2.35%  '
$ws.Range("D4").Value = "'0.9973"
$ws.Range("E4").Value = '  -0.50%  '
$ws.Range("D5").Value = "'307.10"
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("D6").Value = "'0.9962"
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").Value = "'0.3783"
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").Value = "'53.16"
$ws.Range("E8").Value = '  +4.41%  '
$ws.Range("D9").Value = "'0.3656"
$ws.Range("E9").Value = '  +1.56%  '
$ws.Range("D10").Value = "'1.276"
$ws.Range("E10").Value = '  +4.01%  '
$ws.Range("D11").Value = "'0.08184"
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").Value = "'0.9974"
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = "'23.16"
$ws.Range("E13").Value = '  +5.06%  '
$ws.Range("D14").Value = "'6.657"
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("D15").Value = "'7.440"
$ws.Range("E15").Value = '  +2.33%  '
$ws.Range("D16").Value = "'0.00001255"
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = '1.622.90'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = "'94.82"
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").Value = "'0.06934"
$ws.Range("E19").Value = '  +1.88%  '
$ws.Range("E20").Value = '  +2.15%  '
$ws.Range("D21").Value = "'6.585"
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("E23").Value = '  +0.92%  '
$ws.Range("D24").Value = '23.369.35'
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").Value = "'3.158"
$ws.Range("E25").Value = '  +12.05%  '
$ws.Range("D26").Value = "'2.420"
$ws.Range("E26").Value = '  +1.70%  '
$ws.Range("D27").Value = "'21.38"
$ws.Range("E27").Value = '  +2.33%  '
$ws.Range("D28").Value = "'151.00"
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = "'5.278"
$ws.Range("E29").Value = '  +1.28%  '
$ws.Range("D30").Value = "'136.55"
$ws.Range("E30").Value = '  +2.73%  '
$ws.Range("D31").Value = "'2.415"
$ws.Range("E31").Value = '  +2.65%  '
$ws.Range("D32").Value = "'6.884"
$ws.Range("E32").Value = '  +5.16%  '
$ws.Range("D33").Value = '1.799.09'
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").Value = "'0.9710"
$ws.Range("E34").Value = '  +2.68%  '
$ws.Range("D35").Value = "'0.02790"
$ws.Range("E35").Value = '  +4.20%  '
$ws.Range("D36").Value = "'10.44"
$ws.Range("E36").Value = '  +2.73%  '
$ws.Range("D37").Value = "'0.07438"
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("E38").Value = '  +2.88%  '
$ws.Range("D39").Value = "'0.2531"
$ws.Range("E39").Value = '  +2.09%  '
$ws.Range("D40").Value = "'0.08834"
$ws.Range("E40").Value = '  +0.64%  '
$ws.Range("D41").Value = "'1.406"
$ws.Range("E41").Value = '  +4.53%  '
$ws.Range("D42").Value = "'0.7164"
$ws.Range("E42").Value = '  +3.46%  '
$ws.Range("D43").Value = "'12.77"
$ws.Range("E43").Value = '  +5.44%  '
$ws.Range("D44").Value = "'16.12"
$ws.Range("E44").Value = '  +7.59%  '
$ws.Range("D45").Value = "'0.6620"
$ws.Range("E45").Value = '  +2.61%  '
$ws.Range("E46").Value = '  +4.44%  '
$ws.Range("D47").Value = "'4.033"
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("D48").Value = "'0.9954"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = "'0.08019"
$ws.Range("E49").Value = '  +1.63%  '
$ws.Range("D50").Value = "'131.19"
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("E51").Value = '  +0.96%  '
